$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new columns I, J, K, L -----------------------------
# Order matters: shared-string table indices are assigned in first-use order,
# and must land as 20="Column2", 21="Column3", 22="Column4", 29="relabel 2", etc.
$ws.Range("J1").Value = "Column2"
$ws.Range("K1").Value = "Column3"
$ws.Range("L1").Value = "Column4"

# --- New "Chỉnh sửa lab 1 nữa" values (column H) replacing old numbers ------
$ws.Range("H2").Value = "66.76 (c=32)"
$ws.Range("H3").Value = "67.24 (c=31)"
$ws.Range("H11").Value = "69.96 (c=29)"
$ws.Range("H4").Value = "69.68 (c=22)"
$ws.Range("H5").Value = "70.21 (c=27)"
$ws.Range("H6").Value = "69.14 (c=27)"

# --- New "relabel 2" column (column I) header + data ------------------------
$ws.Range("I1").Value = "relabel 2"
$ws.Range("I3").Value = "67.81 (c=33)"
$ws.Range("I2").Value = "67.02 (c=29)"
$ws.Range("I4").Value = "68.64 (c=29)"
$ws.Range("I5").Value = "69.52 (c=32)"
$ws.Range("I6").Value = "68.65 (29)"
$ws.Range("I9").Value = "67.57 (c=33)"
$ws.Range("I10").Value = "68.32 (c=34)"
$ws.Range("I12").Value = "70.88 (c=26)"
$ws.Range("I11").Value = "70.06 (c=33)"
$ws.Range("I13").Value = "69.76 (c=34)"

# H12 reuses the existing "70.54 (c=28)" shared string (previously I12).
$ws.Range("H12").Value = "70.54 (c=28)"
